# "Generate Report for Archive"
#
# 1. Status text "Ready for handoff" -> "In Translation" on all three
#    sheets (Overview!E2:F2, zh-cn!C2, de-de!C2 all shared the same
#    string, so set every occurrence to keep them sharing one entry).
# 2. Narrow the "Status" column(s) from ~17.22 chars to ~13.41 chars:
#      - Overview sheet: columns E and F (zh-cn / de-de status columns)
#      - zh-cn sheet: column C (Status)
#      - de-de sheet: column C (Status)

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# --- 1. Update the status text everywhere it appears ---
$wsOverview.Range("E2").Value = "In Translation"
$wsOverview.Range("F2").Value = "In Translation"
$wsZhCn.Range("C2").Value = "In Translation"
$wsDeDe.Range("C2").Value = "In Translation"

# --- 2. Narrow the status columns ---
# ColumnWidth of 12.5 characters is the closest achievable width to the
# target 13.41 (Excel snaps column widths to its internal pixel grid).
$wsOverview.Columns.Item(5).ColumnWidth = 12.5
$wsOverview.Columns.Item(6).ColumnWidth = 12.5
$wsZhCn.Columns.Item(3).ColumnWidth = 12.5
$wsDeDe.Columns.Item(3).ColumnWidth = 12.5
